$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-63 down to 58-64
$ws.Rows.Item(57).Insert()

# Populate the new row 57 with data (copy of the repeating "Vega Modelo de Temuco" record,
# with its own date / volume / price fields)
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "Vega Modelo de Temuco"
$ws.Range("C57").Value = "La Araucanía"
$ws.Range("D57").Value = 44491
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108004
$ws.Range("J57").Value = "Papaya"
$ws.Range("K57").Value = "Cultivar IV Región"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 40
$ws.Range("N57").Value = 20000
$ws.Range("O57").Value = 20000
$ws.Range("P57").Value = 20000
$ws.Range("Q57").Value = "$/bandeja 10 kilos"
$ws.Range("R57").Value = "Provincia del Elquí"
$ws.Range("S57").Value = 2000
$ws.Range("T57").Value = 10
